$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '42.111.29'
$ws.Range("D2").ClearFormats()
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  -1.67%  '
$ws.Range("E2").ClearFormats()

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.294.18'
$ws.Range("D3").ClearFormats()
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  -2.62%  '
$ws.Range("E3").ClearFormats()

$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  -0.02%  '
$ws.Range("E4").ClearFormats()

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '311.43'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  -3.99%  '
$ws.Range("E5").ClearFormats()

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '105.94'
$ws.Range("D6").ClearFormats()
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  +2.24%  '
$ws.Range("E6").ClearFormats()

$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  -2.04%  '
$ws.Range("E7").ClearFormats()

$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  -0.05%  '
$ws.Range("E8").ClearFormats()

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.608'
$ws.Range("D9").ClearFormats()
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  -2.95%  '
$ws.Range("E9").ClearFormats()

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '40.42'
$ws.Range("D10").ClearFormats()
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  +0.43%  '
$ws.Range("E10").ClearFormats()

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0911'
$ws.Range("D11").ClearFormats()
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  -1.55%  '
$ws.Range("E11").ClearFormats()

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '8.29'
$ws.Range("D12").ClearFormats()
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  -2.66%  '
$ws.Range("E12").ClearFormats()

$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  -0.19%  '
$ws.Range("E13").ClearFormats()

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.968'
$ws.Range("D14").ClearFormats()
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  -3.61%  '
$ws.Range("E14").ClearFormats()

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '15.50'
$ws.Range("D15").ClearFormats()
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  -4.47%  '
$ws.Range("E15").ClearFormats()

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '2.638.64'
$ws.Range("D16").ClearFormats()

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '2.287.59'
$ws.Range("D17").ClearFormats()
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  -2.60%  '
$ws.Range("E17").ClearFormats()

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '42.015.11'
$ws.Range("D18").ClearFormats()
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  -1.74%  '
$ws.Range("E18").ClearFormats()

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '7.55'
$ws.Range("D19").ClearFormats()
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  -4.28%  '
$ws.Range("E19").ClearFormats()

$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  -1.92%  '
$ws.Range("E20").ClearFormats()

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '73.43'
$ws.Range("D21").ClearFormats()
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  -4.74%  '
$ws.Range("E21").ClearFormats()

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '3.47'
$ws.Range("D22").ClearFormats()
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  -5.86%  '
$ws.Range("E22").ClearFormats()

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '256.90'
$ws.Range("D23").ClearFormats()
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  -2.99%  '
$ws.Range("E23").ClearFormats()

$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  -0.61%  '
$ws.Range("E24").ClearFormats()

$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  -7.49%  '
$ws.Range("E25").ClearFormats()

$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  +0.64%  '
$ws.Range("E26").ClearFormats()

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '11.00'
$ws.Range("D27").ClearFormats()
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  -4.45%  '
$ws.Range("E27").ClearFormats()

$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  +2.94%  '
$ws.Range("E28").ClearFormats()

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '22.70'
$ws.Range("D29").ClearFormats()
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  -1.11%  '
$ws.Range("E29").ClearFormats()

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '166.53'
$ws.Range("D30").ClearFormats()
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  -4.88%  '
$ws.Range("E30").ClearFormats()

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '35.57'
$ws.Range("D31").ClearFormats()
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  +0.20%  '
$ws.Range("E31").ClearFormats()

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.0892'
$ws.Range("D32").ClearFormats()
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  -0.94%  '
$ws.Range("E32").ClearFormats()

$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  -6.88%  '
$ws.Range("E33").ClearFormats()

$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  -6.80%  '
$ws.Range("E34").ClearFormats()

$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  +5.91%  '
$ws.Range("E35").ClearFormats()

$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  -2.68%  '
$ws.Range("E36").ClearFormats()

$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  +0.46%  '
$ws.Range("E37").ClearFormats()

$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  -1.80%  '
$ws.Range("E38").ClearFormats()

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.82'
$ws.Range("D39").ClearFormats()
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  +3.28%  '
$ws.Range("E39").ClearFormats()

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '3.61'
$ws.Range("D40").ClearFormats()
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  -4.96%  '
$ws.Range("E40").ClearFormats()

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.49'
$ws.Range("D41").ClearFormats()
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  +0.12%  '
$ws.Range("E41").ClearFormats()

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '71.22'
$ws.Range("D42").ClearFormats()
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  +1.16%  '
$ws.Range("E42").ClearFormats()

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '96.52'
$ws.Range("D43").ClearFormats()
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  +1.73%  '
$ws.Range("E43").ClearFormats()

$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  -4.30%  '
$ws.Range("E44").ClearFormats()

$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  +0.11%  '
$ws.Range("E45").ClearFormats()

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '12.30'
$ws.Range("D46").ClearFormats()
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  +3.02%  '
$ws.Range("E46").ClearFormats()

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '112.59'
$ws.Range("D47").ClearFormats()
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  -6.92%  '
$ws.Range("E47").ClearFormats()

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '9.06'
$ws.Range("D48").ClearFormats()
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  -0.72%  '
$ws.Range("E48").ClearFormats()

$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  -4.87%  '
$ws.Range("E49").ClearFormats()

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '75.01'
$ws.Range("D50").ClearFormats()
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  +5.14%  '
$ws.Range("E50").ClearFormats()

$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  -1.15%  '
$ws.Range("E51").ClearFormats()
